$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    12 = -11.66269999999999
    32 = -13.4067
    36 = -12.5295
    38 = -12.32509999999999
    46 = -14.67789999999999
    54 = -12.76220000000001
    55 = -13.8357
    67 = -11.1128
    69 = -12.09549999999999
    72 = -11.4844
    91 = -10.3661
    99 = -12.9507
}

foreach ($row in $updates.Keys) {
    $ws.Range("C$row").Value = $updates[$row]
}
